$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11 (|S*|/n average)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold font, size 12, vertical center alignment for the B14:B17 cells
$rng = $ws.Range("B14:B17")
$rng.Font.Bold = $true
$rng.Font.Size = 12
$rng.VerticalAlignment = -4108  # xlCenter
$rng.RowHeight = 15.6

$ws.Range("A14:B17").Select()
